# Remove the "Unsubscribe" hyperlink (and its run/text) from the document,
# leaving the now-empty paragraph (with its indentation pPr) in place.
$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("Unsubscribe", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Delete()
}
